$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new shared-string values in the exact order they were first introduced
# (K14, B38, F38, I38, K38, F39, L36, L40, J39, K39, L39) so the underlying
# shared string table is built in the same order as the authored workbook.

$ws.Range("K14").Value = "tableofpartnersandactions"
$ws.Range("B38").Value = "tableofPandA"
$ws.Range("F38").Value = "started as copy of tablesofpartnersandactions.csv so that an extra column could be added and org types could be assigned to each partner row"
$ws.Range("I38").Value = "tableofpartnersandactions.csv (but not created in a script)"
$ws.Range("K38").Value = "finalSdataset_8_29.csv"
$ws.Range("F39").Value = "Includes all actions org types and species in finished product"
$ws.Range("L36").Value = "** I think this is the latest version of salafsky coding actions (came after AS_9_coding) where actions are 1s and 0s"
$ws.Range("L40").Value = "** latest version of salafsky actions categorization before text was changed to 1s and 0s"
$ws.Range("J39").Value = "Partnerships_draft_code"
$ws.Range("K39").Value = "this csv file"
$ws.Range("L39").Value = "combined with correct names of partners and collabs, fixed spelling issues, and assigned correct orgtypes"

# Remaining cells using pre-existing shared strings.
$ws.Range("I14").Value = "codesalafsky.csv"

$ws.Range("A38").Value = "csv file"
$ws.Range("G38").Value = "NA"
$ws.Range("H38").Value = "None"
$ws.Range("J38").Value = "NA"

$ws.Range("A39").Value = "csv file"
$ws.Range("B39").Value = "finalSdataset_8_29.csv"
$ws.Range("I39").Value = "tableofPandAmodified.csv"

$ws.Range("A40").Value = "csv file"
$ws.Range("B40").Value = "codesalafsky.csv"

# --- Update selection to reflect last active cell used while editing ---
$ws.Range("H16").Select()
